$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (column E) values for rows 16-22 to the new
# ascending period order, and the associated "Valor Mora" (column F)
# amounts that travel with the period that used to sit in that slot.
$ws.Range("E16").Value = "2210"
$ws.Range("F16").Value = 80000

$ws.Range("E17").Value = "2211"
$ws.Range("F17").Value = 80000

$ws.Range("E18").Value = "2212"
$ws.Range("F18").Value = 80000

$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 80000

$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 80000

$ws.Range("E21").Value = "2303"
$ws.Range("F21").Value = 80000

$ws.Range("E22").Value = "2304"
$ws.Range("F22").Value = 72000
